$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that were removed from the source data set.
# Delete the later row first (row 28, "SC 92") so the row index for the
# earlier row ("RM 232", row 26) stays valid for the second delete.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# After the deletions the remaining rows have shifted up; re-impute /
# clear individual cells to match the updated (post row-removal) data.

# Clear cells that became missing values
$ws.Range("D3").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Range("F17").ClearContents()
$ws.Range("E18").ClearContents()
$ws.Range("E19").ClearContents()
$ws.Range("F19").ClearContents()
$ws.Range("C27").ClearContents()
$ws.Range("E29").ClearContents()
$ws.Range("F31").ClearContents()
$ws.Range("F32").ClearContents()

# Fill in cells that now carry newly imputed values
$ws.Range("F6").Value = 16.43
$ws.Range("E8").Value = -6.6
$ws.Range("E10").Value = -6.1
$ws.Range("F11").Value = 17.65
$ws.Range("F13").Value = 17.1
$ws.Range("E15").Value = -8.4
$ws.Range("E25").Value = -7.1
$ws.Range("F25").Value = 16.6
$ws.Range("C26").Value = 10.8
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
